# Adding write to excel method and more changes
#
# Adds a small reusable "write to excel" helper and uses it to log a couple
# of test passes/failures next to the existing credentials rows, then
# re-records the original UserName/Password pair as a final confirmation
# row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Write-ToExcel($Sheet, $Row, $Col, $Value) {
    $Sheet.Cells.Item($Row, $Col).Value = $Value
}

# Row 2 (Pradnya / 1994) was correct -> Pass
Write-ToExcel $ws 2 3 "Pass"

# Row 3 (Incorrect / Incorrect) was wrong -> Fail
Write-ToExcel $ws 3 3 "Fail"

# A couple of extra test rows written through the same helper
Write-ToExcel $ws 4 1 "Test"
Write-ToExcel $ws 4 2 "Test"
Write-ToExcel $ws 4 3 "Fail"

Write-ToExcel $ws 5 1 "Test1"
Write-ToExcel $ws 5 2 "Test1"
Write-ToExcel $ws 5 3 "Fail"

# Final row: re-write the original username/password as text (leading
# apostrophe keeps "1994" stored as text, matching the quote-prefixed
# style already used for the row above it) and mark it a Pass.
Write-ToExcel $ws 6 1 "Pradnya"
Write-ToExcel $ws 6 2 "'1994"
Write-ToExcel $ws 6 3 "Pass"

[void]$ws.Range("A5").Select()
